$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# New configuration rows for PDF field extraction (rows 7-12)
$ws.Range("A7").Value = "InvoiceNumber"
$ws.Range("B7").Value = "(?<=Invoice .umber).*?(?=\s*Order .*ber)"

$ws.Range("A8").Value = "OrderNumber"
$ws.Range("B8").Value = "(?<=Order .*ber).*?(?=\s*Invoice .ate)"

$ws.Range("A9").Value = "InvoiceDate"
$ws.Range("B9").Value = "(?<=Invoice .ate).*?(?=\s*Due .ate)"

$ws.Range("A10").Value = "DueDate"
$ws.Range("B10").Value = "(?<=Due .ate).*?(?=\s*Total .mount )"

$ws.Range("A11").Value = "TotalAmount"
$ws.Range("B11").Value = "(?<=Total .mount).*?(?=\s*QTY)"

$ws.Range("A12").Value = "FieldsToExtract"
$ws.Range("B12").Value = "InvoiceNumber,OrderNumber,InvoiceDate,DueDate,TotalAmount"

# Make Settings the active sheet, and focus the B9 cell as the final selection
$ws.Activate()
$ws.Range("B9").Select()
